$wb = $excel.ActiveWorkbook

# The "NewLoanInput" sheet has a test-case identifier string in B2 that
# duplicated the identifier also used elsewhere, creating a test
# inter-dependency. Disambiguate it by appending "-1st".
$ws = $wb.Worksheets.Item("NewLoanInput")
$ws.Range("B2").Value = "4480-RBI-SUBMITLOANON02JAN-INDCOLLSHEETON09JAN-1st"

# Make this sheet the active/selected sheet (as it was prior to the edit).
$ws.Activate()
